# Apply the diff: update the date line and the 25 division problems in the
# single table. Two cells originally share the text "52÷5=" but map to two
# different replacements ("19÷8=" and "41÷8="), so those two cells are
# addressed directly via Table.Cell(row, column) to avoid ambiguity; every
# other string in the document is unique and is updated with Find/Replace.

$d = $word.ActiveDocument

# Header date line.
$d.Content.Find.Execute("2024-07-05 Friday", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-07-06 Saturday", 2)

# Unique division-problem replacements (wdReplaceAll is safe since each
# "find" string occurs exactly once in the document).
$pairs = @(
    @("45÷9=", "57÷8="),
    @("26÷6=", "86÷7="),
    @("13÷9=", "28÷5="),
    @("33÷8=", "37÷5="),
    @("79÷7=", "14÷9="),
    @("45÷6=", "37÷9="),
    @("88÷9=", "36÷2="),
    @("57÷9=", "30÷9="),
    @("61÷3=", "59÷4="),
    @("43÷8=", "52÷2="),
    @("22÷4=", "31÷2="),
    @("19÷7=", "24÷9="),
    @("91÷5=", "87÷9="),
    @("70÷6=", "76÷4="),
    @("16÷7=", "76÷9="),
    @("35÷4=", "62÷4="),
    @("34÷5=", "55÷9="),
    @("82÷3=", "20÷8="),
    @("21÷9=", "52÷6="),
    @("61÷7=", "40÷5="),
    @("67÷7=", "91÷6="),
    @("59÷2=", "21÷4=")
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair[0], $false, $false, $false, $false, $false, `
                             $true, 1, $false, $pair[1], 2)
}

# The two ambiguous "52÷5=" cells: table 1, row 1 col 2 -> "19÷8=",
# table 1, row 9 col 2 -> "41÷8=". Addressed directly by cell so the two
# distinct replacements land on the correct occurrence.
$t = $d.Tables.Item(1)
$t.Cell(1, 2).Range.Text = "19÷8="
$t.Cell(9, 2).Range.Text = "41÷8="
